$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "257.46")
# are preserved verbatim instead of being auto-converted to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "42.877.20"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.216.72"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "257.46"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "76.89"
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "42.20"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "6.99"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "2.545.55"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "14.46"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "2.216.65"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "42.828.34"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("D20").Value = "71.18"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "230.89"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -6.51%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "42.94"
$ws.Range("E26").Value = "  +10.09%  "
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").Value = "3.34"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").Value = "173.09"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "20.46"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "0.0877"
$ws.Range("E33").Value = "  +9.93%  "
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "0.0361"
$ws.Range("E36").Value = "  +7.53%  "
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "4.34"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "12.85"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").Value = "2.12"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +16.97%  "
$ws.Range("D42").Value = "0.201"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "5.30"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").Value = "60.11"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "102.95"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "8.36"
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").Value = "0.467"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0979"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("E51").Value = "  -1.14%  "

# Restore default formatting on column D so no stray number-format styling
# is left behind on the cells.
$dRange.ClearFormats()
